# Generate Report for Handback
#
# This script mutates the localization-status workbook to reflect that the
# handback (localized files coming back from the translators) has now
# happened:
#   - The overall status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" on the Overview sheet.
#   - Each per-language sheet (zh-cn, de-de) gains the "Latest Target File"
#     and "Latest Handback File" values for both rows, a real "Latest
#     Handback DateTime" (replacing the zero date / adding the new one),
#     and a hyperlink on the newly-populated "Latest Target File" cell
#     (column I), matching the existing hyperlink already on column A.
#   - A few columns are widened to fit the newly-populated, longer values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Overview sheet: status text + column widths
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$newStatus = "Handed back: in sync with en-US"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# columns E & F grow to fit the longer status text
$overview.Columns.Item(5).ColumnWidth = 29.16
$overview.Columns.Item(6).ColumnWidth = 29.16

# ---------------------------------------------------------------------
# 2. Per-language sheets: zh-cn and de-de
# ---------------------------------------------------------------------
$file1Name = "462cda04-9234-4fcc-8834-1276f7d60d8d.md"
$file2Name = "d5f011bc-9d14-45a4-a4cf-3eee3d0982a4.md"
$file1Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5b692e91264306c4098838e7e65fece805fb1528/e2e/462cda04-9234-4fcc-8834-1276f7d60d8d.md"
$file2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5b692e91264306c4098838e7e65fece805fb1528/e2e/d5f011bc-9d14-45a4-a4cf-3eee3d0982a4.md"

function Update-LanguageSheet($sheetName, $targetXlf1, $targetXlf2, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Latest Target File (I) / Latest Handback File (J) / Latest Handback
    # DateTime (K) for both data rows.
    $ws.Range("I2").Value = $file1Name
    $ws.Range("J2").Value = $targetXlf1
    $ws.Range("K2").Value = $handbackDateTime

    $ws.Range("I3").Value = $file2Name
    $ws.Range("J3").Value = $targetXlf2
    $ws.Range("K3").Value = $handbackDateTime

    # Rebuild the hyperlinks collection in row-major order (A2, I2, A3, I3)
    # so the new Latest-Target-File hyperlinks (column I) sit next to the
    # existing Source-File-Name hyperlinks (column A).
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $file1Url, "", "", $file1Name)
    $ws.Hyperlinks.Add($ws.Range("I2"), $file1Url, "", "", $file1Name)
    $ws.Hyperlinks.Add($ws.Range("A3"), $file2Url, "", "", $file2Name)
    $ws.Hyperlinks.Add($ws.Range("I3"), $file2Url, "", "", $file2Name)

    # Columns grow to fit the newly-populated, longer values.
    $ws.Columns.Item(3).ColumnWidth = 29.16
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

Update-LanguageSheet "zh-cn" `
    "462cda04-9234-4fcc-8834-1276f7d60d8d.64edee2311166d68656dc7d1b6de7470b2a0d993.zh-cn.xlf" `
    "d5f011bc-9d14-45a4-a4cf-3eee3d0982a4.a734c9a79c45a9eeb1d47dc6ae8930e33196f73c.zh-cn.xlf" `
    "2016-08-21 00:53:53"

Update-LanguageSheet "de-de" `
    "462cda04-9234-4fcc-8834-1276f7d60d8d.64edee2311166d68656dc7d1b6de7470b2a0d993.de-de.xlf" `
    "d5f011bc-9d14-45a4-a4cf-3eee3d0982a4.a734c9a79c45a9eeb1d47dc6ae8930e33196f73c.de-de.xlf" `
    "2016-08-21 00:54:04"
